$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# A1: keep the shared "date" text in use (via D1) before retitling A1 to "Date"
# so both strings survive; D1 becomes a new "date" header cell.
$ws.Range("D1").Value = "date"
$ws.Range("A1").Value = "Date"

# A1 picks up a fresh (General-ish) number-format style distinct from its old one.
$ws.Range("A1").NumberFormat = "General"

# D1 keeps the same custom date-display style already used by column D.
$ws.Range("D1").NumberFormat = "yyyy\-mm\-dd;@"

# --- Data rows: mirror the date values from column A into column D ---
$ws.Range("D2").Value = 43831
$ws.Range("D3").Value = 43832
$ws.Range("D4").Value = 43833

# Column D keeps its existing custom date format.
$ws.Range("D2:D4").NumberFormat = "yyyy\-mm\-dd;@"

# --- Column A (rows 2-8): switch from the built-in date style (numFmtId 14)
# to the workbook's custom yyyy-mm-dd style already used elsewhere (s=1) ---
$ws.Range("A2:A8").NumberFormat = "yyyy\-mm\-dd;@"

# --- Selection moves to D1 ---
$ws.Range("D1").Select()
